$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = 'Última actualización: 11:30:45'
$ws.Range("A3").Value = 'Total filas: 195'
$rows1 = @(
  @('03:45:25','03:47','14_ABASTO',2,'LP1912'),
  @('03:45:25','04:01','81_EL PELIGRO',16,'LP1912'),
  @('03:45:25','04:46','215A_EL PATO',61,'LP1912'),
  @('03:45:25','04:53','11_ETCHEVERRY',68,'LP1912'),
  @('04:56:49','05:13','14_ABASTO',17,'LP1912'),
  @('03:45:25','05:16','17_ROMERO',91,'LP1912'),
  @('04:45:05','05:16','14_ABASTO',31,'LP1912'),
  @('03:45:25','05:22','23_HERNANDEZ',97,'LP1912'),
  @('05:26:08','05:28','14_ABASTO',2,'LP1912'),
  @('03:45:25','05:34','215B_EL PATO',109,'LP1912'),
  @('04:18:02','05:34','14_ABASTO',76,'LP1912'),
  @('04:18:02','05:35','215B_EL PATO',77,'LP1912'),
  @('03:45:25','05:37','14_ABASTO',112,'LP1912'),
  @('04:18:02','05:46','15_ABASTO',88,'LP1912'),
  @('04:45:05','06:04','16_SANTA ANA',79,'LP1912'),
  @('04:18:02','06:05','16_SANTA ANA',107,'LP1912'),
  @('04:56:49','06:11','215A_EL PATO',75,'LP1912'),
  @('04:18:02','06:12','215A_EL PATO',114,'LP1912'),
  @('04:18:02','06:14','225_HARAS DEL SUR',116,'LP1912'),
  @('04:45:05','06:21','26_HERNANDEZ',96,'LP1912'),
  @('06:25:43','06:26','86_EST CHICA-ESC AGRARIA',1,'LP1912'),
  @('04:45:05','06:27','23_HERNANDEZ',102,'LP1912'),
  @('06:25:43','06:28','23_HERNANDEZ',3,'LP1912'),
  @('04:56:49','06:29','86_EST CHICA-ESC AGRARIA',93,'LP1912'),
  @('04:45:05','06:30','86_EST CHICA-ESC AGRARIA',105,'LP1912'),
  @('04:45:05','06:31','16_SANTA ANA',106,'LP1912'),
  @('04:45:05','06:44','225_C ROCA-H SUR',119,'LP1912'),
  @('05:55:25','06:44','26_HERNANDEZ',49,'LP1912'),
  @('04:56:49','06:46','215C_EL PATO',110,'LP1912'),
  @('05:26:08','06:47','215C_EL PATO',81,'LP1912'),
  @('05:55:25','06:59','14_ABASTO',64,'LP1912'),
  @('05:26:08','07:00','14_ABASTO',94,'LP1912'),
  @('06:25:43','07:01','16_SANTA ANA',36,'LP1912'),
  @('05:55:25','07:04','23_HERNANDEZ',69,'LP1912'),
  @('05:26:08','07:05','23_HERNANDEZ',99,'LP1912'),
  @('05:26:08','07:05','15_ABASTO',99,'LP1912'),
  @('05:26:08','07:06','10_OLMOS',100,'LP1912'),
  @('05:26:08','07:07','225_GOMEZ',101,'LP1912'),
  @('05:26:08','07:11','215A_EL PATO',105,'LP1912'),
  @('06:55:02','07:12','215A_EL PATO',17,'LP1912'),
  @('06:25:43','07:14','26_HERNANDEZ',49,'LP1912'),
  @('05:55:25','07:15','11_ETCHEVERRY',80,'LP1912'),
  @('05:26:08','07:16','11_ETCHEVERRY',110,'LP1912'),
  @('06:55:02','07:17','16_SANTA ANA',22,'LP1912'),
  @('05:26:08','07:21','26_HERNANDEZ',115,'LP1912'),
  @('05:26:08','07:23','10_OLMOS',117,'LP1912'),
  @('05:55:25','07:30','10_OLMOS',95,'LP1912'),
  @('05:55:25','07:31','11_ETCHEVERRY',96,'LP1912'),
  @('05:55:25','07:31','16_SANTA ANA',96,'LP1912'),
  @('05:55:25','07:32','84_COLONIA URQUIZA-ESC 49',97,'LP1912'),
  @('06:55:02','07:32','16_SANTA ANA',37,'LP1912'),
  @('06:55:02','07:32','11_ETCHEVERRY',37,'LP1912'),
  @('07:19:29','07:35','23_HERNANDEZ',16,'LP1912'),
  @('05:55:25','07:36','27_EL RETIRO',101,'LP1912'),
  @('06:55:02','07:37','27_EL RETIRO',42,'LP1912'),
  @('05:55:25','07:39','10_OLMOS',104,'LP1912'),
  @('07:19:29','07:46','16_SANTA ANA',27,'LP1912'),
  @('05:55:25','07:47','14_ABASTO',112,'LP1912'),
  @('06:55:02','07:48','14_ABASTO',53,'LP1912'),
  @('07:50:16','07:50','10_OLMOS',0,'LP1912'),
  @('05:55:25','07:51','215D_EL PATO',116,'LP1912'),
  @('06:55:02','07:52','215D_EL PATO',57,'LP1912'),
  @('07:19:29','07:59','23_HERNANDEZ',40,'LP1912'),
  @('06:25:43','08:01','23_HERNANDEZ',96,'LP1912'),
  @('07:19:29','08:03','11_ETCHEVERRY',44,'LP1912'),
  @('06:55:02','08:03','23_HERNANDEZ',68,'LP1912'),
  @('08:02:22','08:05','23_HERNANDEZ',3,'LP1912'),
  @('07:19:29','08:10','16_SANTA ANA',51,'LP1912'),
  @('07:50:16','08:11','16_SANTA ANA',21,'LP1912'),
  @('06:25:43','08:12','15_ABASTO',107,'LP1912'),
  @('07:50:16','08:13','10_OLMOS',23,'LP1912'),
  @('06:55:02','08:21','26_HERNANDEZ',86,'LP1912'),
  @('06:25:43','08:22','16_P MOR-SANTA ANA',117,'LP1912'),
  @('06:55:02','08:23','16_P MOR-SANTA ANA',88,'LP1912'),
  @('06:25:43','08:23','215B_EL PATO',118,'LP1912'),
  @('06:55:02','08:27','84_COLONIA URQUIZA-ESC 49',92,'LP1912'),
  @('07:50:16','08:30','23_HERNANDEZ',40,'LP1912'),
  @('08:02:22','08:33','10_OLMOS',31,'LP1912'),
  @('08:02:22','08:34','23_HERNANDEZ',32,'LP1912'),
  @('08:32:09','08:37','23_HERNANDEZ',5,'LP1912'),
  @('06:55:02','08:42','81_EL PELIGRO',107,'LP1912'),
  @('07:19:29','08:43','14_ABASTO',84,'LP1912'),
  @('07:50:16','08:44','14_ABASTO',54,'LP1912'),
  @('08:32:09','08:53','10_OLMOS',21,'LP1912'),
  @('06:55:02','08:54','17_ROMERO',119,'LP1912'),
  @('07:19:29','09:01','215A_EL PATO',102,'LP1912'),
  @('07:50:16','09:02','215A_EL PATO',72,'LP1912'),
  @('08:02:22','09:03','11_ETCHEVERRY',61,'LP1912'),
  @('08:32:09','09:04','11_ETCHEVERRY',32,'LP1912'),
  @('08:32:09','09:05','23_HERNANDEZ',33,'LP1912'),
  @('07:19:29','09:10','16_P MOR-SANTA ANA',111,'LP1912'),
  @('07:50:16','09:11','16_P MOR-SANTA ANA',81,'LP1912'),
  @('08:32:09','09:13','10_OLMOS',41,'LP1912'),
  @('07:19:29','09:16','27_EL RETIRO',117,'LP1912'),
  @('07:50:16','09:17','27_EL RETIRO',87,'LP1912'),
  @('07:50:16','09:21','26_HERNANDEZ',91,'LP1912'),
  @('08:02:22','09:22','16_SANTA ANA',80,'LP1912'),
  @('07:50:16','09:23','17_ROMERO',93,'LP1912'),
  @('08:02:22','09:23','11_ETCHEVERRY',81,'LP1912'),
  @('08:32:09','09:23','16_SANTA ANA',51,'LP1912'),
  @('07:50:16','09:24','11_ETCHEVERRY',94,'LP1912'),
  @('07:50:16','09:28','16_SANTA ANA',98,'LP1912'),
  @('07:50:16','09:32','15_ABASTO',102,'LP1912'),
  @('07:50:16','09:33','10_OLMOS',103,'LP1912'),
  @('08:56:29','09:34','23_HERNANDEZ',38,'LP1912'),
  @('08:56:29','09:34','16_SANTA ANA',38,'LP1912'),
  @('08:32:09','09:35','16_SANTA ANA',63,'LP1912'),
  @('08:48:08','09:35','23_HERNANDEZ',47,'LP1912'),
  @('09:35:26','09:39','23_HERNANDEZ',4,'LP1912'),
  @('07:50:16','09:42','215C_EL PATO',112,'LP1912'),
  @('08:02:22','09:43','14_ABASTO',101,'LP1912'),
  @('07:50:16','09:44','14_ABASTO',114,'LP1912'),
  @('09:35:26','09:46','16_SANTA ANA',11,'LP1912'),
  @('08:32:09','09:52','15_ABASTO',80,'LP1912'),
  @('08:56:29','09:53','10_OLMOS',57,'LP1912'),
  @('09:35:26','09:58','16_SANTA ANA',23,'LP1912'),
  @('09:35:26','10:03','11_ETCHEVERRY',28,'LP1912'),
  @('08:56:29','10:10','16_P MOR-SANTA ANA',74,'LP1912'),
  @('08:32:09','10:11','16_P MOR-SANTA ANA',99,'LP1912'),
  @('09:35:26','10:12','15_ABASTO',37,'LP1912'),
  @('09:35:26','10:13','10_OLMOS',38,'LP1912'),
  @('08:32:09','10:21','26_HERNANDEZ',109,'LP1912'),
  @('08:32:09','10:22','17_ROMERO',110,'LP1912'),
  @('09:35:26','10:23','11_ETCHEVERRY',48,'LP1912'),
  @('08:56:29','10:26','215A_EL PATO',90,'LP1912'),
  @('08:32:09','10:27','215A_EL PATO',115,'LP1912'),
  @('10:29:57','10:29','16_SANTA ANA',0,'LP1912'),
  @('10:29:57','10:31','10_OLMOS',2,'LP1912'),
  @('09:35:26','10:34','23_HERNANDEZ',59,'LP1912'),
  @('10:29:57','10:34','16_SANTA ANA',5,'LP1912'),
  @('10:29:57','10:39','23_HERNANDEZ',10,'LP1912'),
  @('10:29:57','10:41','17_ROMERO',12,'LP1912'),
  @('08:48:08','10:42','17_ROMERO',114,'LP1912'),
  @('08:56:29','10:43','14_ABASTO',107,'LP1912'),
  @('08:48:08','10:44','14_ABASTO',116,'LP1912'),
  @('10:29:57','10:51','15_ABASTO',22,'LP1912'),
  @('10:29:57','10:52','10_OLMOS',23,'LP1912'),
  @('09:35:26','10:54','27_EL RETIRO',79,'LP1912'),
  @('10:29:57','10:56','27_EL RETIRO',27,'LP1912'),
  @('10:59:49','10:59','16_SANTA ANA',0,'LP1912'),
  @('10:29:57','11:01','215C_EL PATO',32,'LP1912'),
  @('09:35:26','11:02','215C_EL PATO',87,'LP1912'),
  @('10:29:57','11:03','11_ETCHEVERRY',34,'LP1912'),
  @('10:29:57','11:04','23_HERNANDEZ',35,'LP1912'),
  @('09:35:26','11:06','16_P MOR-167 Y 521',91,'LP1912'),
  @('10:59:49','11:06','23_HERNANDEZ',7,'LP1912'),
  @('10:59:49','11:11','10_OLMOS',12,'LP1912'),
  @('10:29:57','11:11','15_ABASTO',42,'LP1912'),
  @('10:59:49','11:12','15_ABASTO',13,'LP1912'),
  @('09:35:26','11:19','86_EST CHICA-ESC AGRARIA',104,'LP1912'),
  @('10:29:57','11:20','26_HERNANDEZ',51,'LP1912'),
  @('09:35:26','11:21','26_HERNANDEZ',106,'LP1912'),
  @('10:29:57','11:26','225_C ROCA-H SUR',57,'LP1912'),
  @('09:35:26','11:27','225_C ROCA-H SUR',112,'LP1912'),
  @('11:30:45','11:30','16_SANTA ANA',0,'LP1912'),
  @('11:30:45','11:31','16_SANTA ANA',1,'LP1912'),
  @('10:29:57','11:31','81_EL PELIGRO',62,'LP1912'),
  @('09:35:26','11:32','81_EL PELIGRO',117,'LP1912'),
  @('10:59:49','11:34','23_HERNANDEZ',35,'LP1912'),
  @('10:29:57','11:35','11_ETCHEVERRY',66,'LP1912'),
  @('10:29:57','11:40','10_OLMOS',71,'LP1912'),
  @('10:29:57','11:41','17_ROMERO',72,'LP1912'),
  @('10:59:49','11:42','11_ETCHEVERRY',43,'LP1912'),
  @('10:59:49','11:43','10_OLMOS',44,'LP1912'),
  @('11:30:45','11:44','11_ETCHEVERRY',14,'LP1912'),
  @('10:29:57','11:50','215B_EL PATO',81,'LP1912'),
  @('10:59:49','11:51','215B_EL PATO',52,'LP1912'),
  @('10:59:49','11:52','15_ABASTO',53,'LP1912'),
  @('10:29:57','11:58','225_GOMEZ',89,'LP1912'),
  @('11:30:45','11:59','225_GOMEZ',29,'LP1912'),
  @('10:29:57','12:01','84_COLONIA URQUIZA-ESC 49',92,'LP1912'),
  @('10:59:49','12:02','84_COLONIA URQUIZA-ESC 49',63,'LP1912'),
  @('11:30:45','12:04','23_HERNANDEZ',34,'LP1912'),
  @('10:59:49','12:06','14_ABASTO',67,'LP1912'),
  @('10:29:57','12:06','16_P MOR-SANTA ANA',97,'LP1912'),
  @('10:59:49','12:10','10_OLMOS',71,'LP1912'),
  @('10:29:57','12:13','17_ROMERO',104,'LP1912'),
  @('10:29:57','12:15','14_ABASTO',106,'LP1912'),
  @('10:29:57','12:20','26_HERNANDEZ',111,'LP1912'),
  @('10:29:57','12:20','215A_EL PATO',111,'LP1912'),
  @('10:59:49','12:20','14_ABASTO',81,'LP1912'),
  @('10:59:49','12:21','26_HERNANDEZ',82,'LP1912'),
  @('10:59:49','12:30','17_ROMERO',91,'LP1912'),
  @('10:59:49','12:36','27_EL RETIRO',97,'LP1912'),
  @('10:59:49','12:38','17_179 Y 38',99,'LP1912'),
  @('11:30:45','12:41','10_OLMOS',71,'LP1912'),
  @('10:59:49','12:48','11_ETCHEVERRY',109,'LP1912'),
  @('11:30:45','13:01','17_ROMERO',91,'LP1912'),
  @('11:30:45','13:06','16_P MOR-SANTA ANA',96,'LP1912'),
  @('11:30:45','13:07','10_OLMOS',97,'LP1912'),
  @('11:30:45','13:13','215D_EL PATO',103,'LP1912'),
  @('11:30:45','13:21','26_HERNANDEZ',111,'LP1912'),
  @('11:30:45','13:25','10_OLMOS',115,'LP1912'),
  @('11:30:45','13:26','14_ABASTO',116,'LP1912'),
  @('11:30:45','13:26','15_ABASTO',116,'LP1912')
)
$r = 6
foreach ($row in $rows1) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $r = $r + 1
}

$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = 'Última actualización: 11:30:45'
$ws.Range("A3").Value = 'Total filas: 23'
$rows2 = @(
  @('03:45:25','04:46','215A_EL PATO',61,'LP1912'),
  @('03:45:25','05:34','215B_EL PATO',109,'LP1912'),
  @('04:18:02','05:35','215B_EL PATO',77,'LP1912'),
  @('04:56:49','06:11','215A_EL PATO',75,'LP1912'),
  @('04:18:02','06:12','215A_EL PATO',114,'LP1912'),
  @('04:56:49','06:46','215C_EL PATO',110,'LP1912'),
  @('05:26:08','06:47','215C_EL PATO',81,'LP1912'),
  @('05:26:08','07:11','215A_EL PATO',105,'LP1912'),
  @('06:55:02','07:12','215A_EL PATO',17,'LP1912'),
  @('05:55:25','07:51','215D_EL PATO',116,'LP1912'),
  @('06:55:02','07:52','215D_EL PATO',57,'LP1912'),
  @('06:25:43','08:23','215B_EL PATO',118,'LP1912'),
  @('07:19:29','09:01','215A_EL PATO',102,'LP1912'),
  @('07:50:16','09:02','215A_EL PATO',72,'LP1912'),
  @('07:50:16','09:42','215C_EL PATO',112,'LP1912'),
  @('08:56:29','10:26','215A_EL PATO',90,'LP1912'),
  @('08:32:09','10:27','215A_EL PATO',115,'LP1912'),
  @('10:29:57','11:01','215C_EL PATO',32,'LP1912'),
  @('09:35:26','11:02','215C_EL PATO',87,'LP1912'),
  @('10:29:57','11:50','215B_EL PATO',81,'LP1912'),
  @('10:59:49','11:51','215B_EL PATO',52,'LP1912'),
  @('10:29:57','12:20','215A_EL PATO',111,'LP1912'),
  @('11:30:45','13:13','215D_EL PATO',103,'LP1912')
)
$r = 6
foreach ($row in $rows2) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $r = $r + 1
}

$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = 'Última actualización: 11:30:45'
$ws.Range("A3").Value = 'Total filas: 32'
$rows3 = @(
  @('04:56:49','05:43','215A_LA PLATA',47,'L6173'),
  @('03:45:25','05:44','215A_LA PLATA',119,'L6173'),
  @('04:56:49','06:08','215A_LA PLATA',72,'L6173'),
  @('04:18:02','06:09','215A_LA PLATA',111,'L6173'),
  @('04:56:49','06:32','215C_LA PLATA',96,'L6203'),
  @('04:45:05','06:33','215C_LA PLATA',108,'L6203'),
  @('06:25:43','06:59','215B_LP-P MOR-1 Y 57',34,'L6173'),
  @('05:26:08','07:00','215B_LP-P MOR-1 Y 57',94,'L6173'),
  @('05:55:25','07:35','215A_LA PLATA',100,'L6173'),
  @('06:25:43','07:39','215A_LA PLATA',74,'L6173'),
  @('06:55:02','07:42','215A_LA PLATA',47,'L6173'),
  @('07:19:29','07:46','215A_LA PLATA',27,'L6173'),
  @('07:50:16','07:51','215A_LA PLATA',1,'L6173'),
  @('06:25:43','08:06','215C_LA PLATA',101,'L6203'),
  @('06:55:02','08:07','215C_LA PLATA',72,'L6203'),
  @('07:19:29','08:21','215C_LA PLATA',62,'L6203'),
  @('08:02:22','08:25','215C_LA PLATA',23,'L6203'),
  @('07:50:16','08:27','215C_LA PLATA',37,'L6203'),
  @('07:19:29','08:35','215A_LA PLATA',76,'L6173'),
  @('06:55:02','08:36','215A_LA PLATA',101,'L6173'),
  @('08:02:22','08:39','215A_LA PLATA',37,'L6173'),
  @('08:32:09','08:42','215A_LA PLATA',10,'L6173'),
  @('07:19:29','09:09','215D_LA PLATA',110,'L6203'),
  @('08:32:09','10:03','215B_LP-P MOR-40 Y 115',91,'L6173'),
  @('10:29:57','10:53','215A_LA PLATA',24,'L6173'),
  @('08:56:29','10:54','215A_LA PLATA',118,'L6173'),
  @('10:29:57','11:13','215C_LA PLATA',44,'L6203'),
  @('09:35:26','11:14','215C_LA PLATA',99,'L6203'),
  @('10:29:57','12:03','215A_LA PLATA',94,'L6173'),
  @('10:59:49','12:07','215A_LA PLATA',68,'L6173'),
  @('11:30:45','12:09','215A_LA PLATA',39,'L6173'),
  @('10:59:49','12:53','215C_LA PLATA',114,'L6203')
)
$r = 6
foreach ($row in $rows3) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $r = $r + 1
}

Write-Host "done"